$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Default_Values")
$ws.Activate()

$ws.Range("B3").Value = 99999999999
$ws.Range("B23").Value = 99999999999
$ws.Range("B43").Value = 99999999999
$ws.Range("B44").Value = 99999999999
$ws.Range("B48").Value = 99999999999
$ws.Range("B50").Value = 99999999999

$ws.Range("B3").Select()
